# Weekly update: a new week of "Ají" (Comercializadora del Agro de Limarí)
# price data is prepended at the top of the historical log (rows 140-141),
# pushing the previously-newest rows (old 140-162) down by two rows
# (to 142-164).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 140:141 - shifts old rows 140-162 down to 142-164.
$ws.Range("A140:R141").Insert()

# Row 140 - Ají / Americana (o) / Primera, new week (date 44476)
$ws.Range("A140").Value = 2
$ws.Range("B140").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C140").Value = "Coquimbo"
$ws.Range("D140").Value = 44476
$ws.Range("E140").Value = 4
$ws.Range("F140").Value = 100112021
$ws.Range("G140").Value = "Ají"
$ws.Range("H140").Value = "Americana (o)"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 160
$ws.Range("K140").Value = 80000
$ws.Range("L140").Value = 85000
$ws.Range("M140").Value = 82500
$ws.Range("N140").Value = "$/caja 25 kilos"
$ws.Range("O140").Value = "Provincia de Limarí"
$ws.Range("P140").Value = 3300
$ws.Range("Q140").Value = 25
$ws.Range("R140").Value = "Hortaliza"

# Row 141 - Ají / Inferno / Primera, new week (date 44476)
$ws.Range("A141").Value = 2
$ws.Range("B141").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 44476
$ws.Range("E141").Value = 4
$ws.Range("F141").Value = 100112021
$ws.Range("G141").Value = "Ají"
$ws.Range("H141").Value = "Inferno"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 140
$ws.Range("K141").Value = 65000
$ws.Range("L141").Value = 70000
$ws.Range("M141").Value = 67500
$ws.Range("N141").Value = "$/caja 25 kilos"
$ws.Range("O141").Value = "Provincia de Limarí"
$ws.Range("P141").Value = 2700
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"
